$d = $word.ActiveDocument

$summaryIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq "Summary:") {
        $summaryIndex = $i
        break
    }
}

if ($summaryIndex -eq 0) {
    throw "Summary: paragraph not found"
}

$headingPara = $d.Paragraphs.Item($summaryIndex)
$bodyPara = $d.Paragraphs.Item($summaryIndex + 1)

$headingXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="GDD"/></w:pPr><w:bookmarkStart w:id="1" w:name="_Toc377218695"/><w:r><w:t>Summary:</w:t></w:r><w:bookmarkEnd w:id="1"/></w:p>
'@

$bodyXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:t xml:space="preserve">The </w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">mighty </w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>princesses of the Fox Clan are sent out</w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> to rescue the clan’s princes who have been</w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> captured by the enemy Dog Clan. </w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">The player needs to successfully help the princesses hack and slash their way through enemy clan members to reach the end of each level. </w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Each section of levels will feature a boss fight where the reward is one of the captured princes. </w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Players will be </w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>able to customize the</w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> princess stats, weapons, and pets as they gain abilities, money and stat points as the game progresses. Players also have the opt</w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">ion to spend money they find </w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">along </w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>the way in the game to buy new items such as weapons and pets</w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> or buffs. </w:t></w:r><w:bookmarkStart w:id="2" w:name="_GoBack"/><w:bookmarkEnd w:id="2"/></w:p>
'@

$r1 = $headingPara.Range
$null = $r1.InsertXML($headingXml)

$r2 = $bodyPara.Range
$null = $r2.InsertXML($bodyXml)

Write-Output "done"
